$d = $word.ActiveDocument

# --- Locate the last paragraph in the document (the one that currently
#     ends with "...Regresion logistica -> ... ?Lo pongo?" and holds the
#     _GoBack bookmark) and append a brand new list paragraph after it. ---
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphAfter()

# The freshly created paragraph inherits the list style/numbering (Prrafodelista,
# numId 3) automatically from the paragraph it was split off from.
$newIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newIndex)
$newStart = $newPara.Range.Start

$text1 = "Preguntar, aparte de los campos isbn, doi y keywords, si hay que quitar los siguientes campos: address, publisher, location, series, url y"
$text2 = " acmid"

# Insert the full text of the new paragraph first (this also moves the
# document's real end further along, so that the bookmark we re-add below
# will no longer sit at the very end of the story -- doing that earlier
# triggers a relocation quirk with the special "_GoBack" bookmark).
$fillRange = $d.Range($newStart, $newStart)
$fillRange.InsertAfter($text1 + $text2)

# Position right between the two runs, where the _GoBack bookmark must live.
$bookmarkPos = $newStart + $text1.Length

# Move the _GoBack bookmark from the end of the previous paragraph to this
# new position in the middle of the freshly typed sentence.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
